# Add a new company/domain row (Incap Co / www.thepacket.ninja) as row 7,
# matching the style of the existing rows, and wire up a hyperlink for the
# new domain cell (mirrors the existing B2:B6 hyperlink cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row
$ws.Range("A7").Value = "Incap Co"
$ws.Range("B7").Value = "www.thepacket.ninja"

# Match A7's formatting to the other "Company Name" cells in column A
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)

# Turn the new domain cell into a live hyperlink, like B2:B6
$ws.Hyperlinks.Add($ws.Range("B7"), "http://www.thepacket.ninja")
$hl = $ws.Hyperlinks.Item($ws.Hyperlinks.Count)
$hl.Address = "http://www.thepacket.ninja"
